$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto price/volume refresh (GitHub Actions scheduled update).
# Source feed values are plain text snapshots (Price/Volume columns), so
# any cell whose new value looks like a plain number is explicitly set to
# Text format first -- otherwise Excel would silently re-parse it as a
# number (dropping significant trailing zeros / introducing float noise).

$ws.Range("D2").Value = "65.830.47"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.661.22"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.71"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.05"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "3.139.73"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "65.714.98"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "2.672.83"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  -2.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.80"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.47"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.05"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +9.66%  "
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.65"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "581.44"
$ws.Range("E28").Value = "  +9.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.15"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.74"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.422"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.59"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.38"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("E41").Value = "  +9.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.81"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.10"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.34"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0260"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.643"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.79"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -6.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.815"
$ws.Range("E51").Value = "  +0.33%  "
